# Update countries & provincias Spain
# - Re-sort Uganda so it appears right after Sri Lanka (pushing Gambia,
#   Lituania and Eslovenia down one row); Mali stays put.
# - Refresh the covid case counters for a batch of countries.
# - Bump the "Datos actualizados..." timestamp string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "last updated" timestamp string (cell A1) ---------------
$ws.Range("A1").Value = "Datos actualizados a 30 de Agosto de 2020 a las 20:54"

# --- 2. Re-order the Uganda / Gambia / Lituania / Eslovenia block ----------
# Before: 127 Sri Lanka, 128 Gambia, 129 Lituania, 130 Eslovenia, 131 Uganda, 132 Mali
# After : 127 Sri Lanka, 128 Uganda, 129 Gambia, 130 Lituania, 131 Eslovenia, 132 Mali
# Capture the current (pre-reorder) rows for Gambia/Lituania/Eslovenia first,
# then shift each one down a row so Uganda's old slot can become "Eslovenia"
# and Uganda itself can be written into row 128 afterwards.
$gambia    = @($ws.Cells.Item(128,1).Value2, $ws.Cells.Item(128,2).Value2, $ws.Cells.Item(128,3).Value2, $ws.Cells.Item(128,4).Value2, $ws.Cells.Item(128,5).Value2, $ws.Cells.Item(128,6).Value2, $ws.Cells.Item(128,7).Value2, $ws.Cells.Item(128,8).Value2)
$lituania  = @($ws.Cells.Item(129,1).Value2, $ws.Cells.Item(129,2).Value2, $ws.Cells.Item(129,3).Value2, $ws.Cells.Item(129,4).Value2, $ws.Cells.Item(129,5).Value2, $ws.Cells.Item(129,6).Value2, $ws.Cells.Item(129,7).Value2, $ws.Cells.Item(129,8).Value2)
$eslovenia = @($ws.Cells.Item(130,1).Value2, $ws.Cells.Item(130,2).Value2, $ws.Cells.Item(130,3).Value2, $ws.Cells.Item(130,4).Value2, $ws.Cells.Item(130,5).Value2, $ws.Cells.Item(130,6).Value2, $ws.Cells.Item(130,7).Value2, $ws.Cells.Item(130,8).Value2)

function Write-RowArray($row, $values) {
    for ($i = 0; $i -lt 8; $i++) {
        $ws.Cells.Item($row, $i + 1).Value = $values[$i]
    }
}

Write-RowArray 129 $gambia
Write-RowArray 130 $lituania
Write-RowArray 131 $eslovenia

$ws.Cells.Item(128, 1).Value = "Uganda"

# --- 3. Refresh numeric data (Casos totales, Nuevos casos, Casos activos,
#        Recuperados, Casos criticos, Muertes hoy, Muertes) ----------------

function Set-Row($row, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

Set-Row 4   6158440 19362 3414253 2557160 0 172 187027    # Estados Unidos
Set-Row 6   3619169 79457 2772928 781624  0 960 64617     # India
Set-Row 13  409974  1965  382584  16146   0 63  11244     # Chile
Set-Row 23  243092  267   217484  16245   0 0   9363      # Alemania
Set-Row 24  231177  3731  172880  51338   0 68  6959      # Irak
Set-Row 27  127940  267   113663  5160    0 4   9117      # Canada
Set-Row 103 7365    249   2929    4364    0 3   72        # Namibia
Set-Row 127 3012    17    2860    140     0 0   12        # Sri Lanka
Set-Row 128 2928    81    1288    1610    0 1   30        # Uganda (new row position)
Set-Row 129 2895    0     751     2048    0 0   96        # Gambia
Set-Row 130 2874    35    1837    951     0 0   86        # Lituania
Set-Row 131 2865    31    2283    449     0 0   133       # Eslovenia
Set-Row 133 2703    75    614     1980    0 3   109       # Siria
Set-Row 143 1997    22    765     1222    0 0   10        # Aruba
Set-Row 145 1953    7     1123    266     0 1   564       # Yemen
Set-Row 164 1012    4     878     57      0 0   77        # Republica del Chad
